$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: epochs 8 -> epochs 10 params; A 0 -> 1; C/D updated; E stays 1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "{'anOptimizer': 'adam', 'batch_size': 10, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("C2").Value = 0.923655370871226
$ws.Range("D2").Value = 0.01501205741479679
$ws.Range("E2").Value = 1

# Row 3: epochs 10 -> epochs 8 params; A 1 -> 0; C/D updated; E 1 -> 2
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "{'anOptimizer': 'adam', 'batch_size': 10, 'epochs': 8, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("C3").Value = 0.923613707224528
$ws.Range("D3").Value = 0.02000352744147469
$ws.Range("E3").Value = 2

# Row 4: batch_size 50 epochs 8 -> epochs 10 params; A 2 -> 3; C/D updated; E stays 3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "{'anOptimizer': 'adam', 'batch_size': 50, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("C4").Value = 0.7840478221575419
$ws.Range("D4").Value = 0.07511020857037724
$ws.Range("E4").Value = 3

# Row 5: batch_size 50 epochs 10 -> epochs 8 params; A 3 -> 2; C/D updated; E stays 4
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "{'anOptimizer': 'adam', 'batch_size': 50, 'epochs': 8, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("C5").Value = 0.7159105141957601
$ws.Range("D5").Value = 0.06637576653389871
$ws.Range("E5").Value = 4
